# Auto-generated edit script applying the Chocobo_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across
# multiple crafting-job sheets per the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 1850.5
$ws.Range("I86").Value = 1601.5
$ws.Range("J86").Value = 1975
$ws.Range("K86").Value = 1601.5
$ws.Range("L86").Value = 1975
$ws.Range("M86").Value = -478.5
$ws.Range("N86").Value = -4221
# Row 89
$ws.Range("H89").Value = 1850.5
$ws.Range("I89").Value = 1601.5
$ws.Range("J89").Value = 1975
$ws.Range("K89").Value = 8007.5
$ws.Range("L89").Value = 9875
$ws.Range("M89").Value = -2391.5
$ws.Range("N89").Value = -21107
# Row 98
$ws.Range("H98").Value = 11750
$ws.Range("I98").Value = 10000
$ws.Range("J98").Value = 13500
$ws.Range("K98").Value = 10000
$ws.Range("L98").Value = 13500
$ws.Range("M98").Value = -8502
$ws.Range("N98").Value = -16496
# Row 107
$ws.Range("H107").Value = 1970.1333
$ws.Range("I107").Value = 2972.8572
$ws.Range("J107").Value = 1092.75
$ws.Range("K107").Value = 2972.8572
$ws.Range("L107").Value = 1092.75
$ws.Range("M107").Value = -1052.8572
$ws.Range("N107").Value = -4932.75
# Row 112
$ws.Range("H112").Value = 1288.035
$ws.Range("I112").Value = 600
$ws.Range("J112").Value = 1326.2593
$ws.Range("K112").Value = 1800
$ws.Range("L112").Value = 3978.7779
$ws.Range("M112").Value = -692
$ws.Range("N112").Value = -6194.7779
# Row 122
$ws.Range("H122").Value = 11750
$ws.Range("I122").Value = 10000
$ws.Range("J122").Value = 13500
$ws.Range("K122").Value = 30000
$ws.Range("L122").Value = 40500
$ws.Range("M122").Value = -27550
$ws.Range("N122").Value = -45400
# Row 137
$ws.Range("H137").Value = 1402721.9
$ws.Range("I137").Value = 1702169.4
$ws.Range("J137").Value = 5300
$ws.Range("K137").Value = 5106508.199999999
$ws.Range("L137").Value = 15900
$ws.Range("M137").Value = -5103958.199999999
$ws.Range("N137").Value = -21000

$ws = $wb.Worksheets.Item("ARM")
# Row 68
$ws.Range("H68").Value = 15500
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1189
# Row 71
$ws.Range("H71").Value = 15500
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 6000
$ws.Range("M71").Value = -1944
# Row 74
$ws.Range("H74").Value = 2193.3845
$ws.Range("I74").Value = 801.7143
$ws.Range("K74").Value = 801.7143
$ws.Range("M74").Value = 72.28570000000002
# Row 77
$ws.Range("H77").Value = 2193.3845
$ws.Range("I77").Value = 801.7143
$ws.Range("K77").Value = 4008.5715
$ws.Range("M77").Value = 359.4285
# Row 132
$ws.Range("H132").Value = 1479.1428
$ws.Range("I132").Value = 859.25
$ws.Range("J132").Value = 3462.8
$ws.Range("K132").Value = 2577.75
$ws.Range("L132").Value = 10388.4
$ws.Range("M132").Value = -47.75
$ws.Range("N132").Value = -15448.4
# Row 139
$ws.Range("H139").Value = 41649.5
$ws.Range("J139").Value = 41649.5
$ws.Range("L139").Value = 41649.5
$ws.Range("N139").Value = -51929.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 13024.429
$ws.Range("I20").Value = 2653
$ws.Range("J20").Value = 20803
$ws.Range("K20").Value = 2653
$ws.Range("L20").Value = 20803
$ws.Range("M20").Value = -2406
$ws.Range("N20").Value = -21297
# Row 62
$ws.Range("H62").Value = 42000
$ws.Range("J62").Value = 42000
$ws.Range("L62").Value = 42000
$ws.Range("N62").Value = -43372
# Row 65
$ws.Range("H65").Value = 42000
$ws.Range("J65").Value = 42000
$ws.Range("L65").Value = 126000
$ws.Range("N65").Value = -132864
# Row 86
$ws.Range("H86").Value = 2633.3333
$ws.Range("I86").Value = 1800
$ws.Range("J86").Value = 2871.4285
$ws.Range("K86").Value = 1800
$ws.Range("L86").Value = 2871.4285
$ws.Range("M86").Value = -677
$ws.Range("N86").Value = -5117.4285
# Row 89
$ws.Range("H89").Value = 2633.3333
$ws.Range("I89").Value = 1800
$ws.Range("J89").Value = 2871.4285
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 14357.1425
$ws.Range("M89").Value = -3384
$ws.Range("N89").Value = -25589.1425
# Row 107
$ws.Range("H107").Value = 3236
$ws.Range("I107").Value = 3011
$ws.Range("J107").Value = 3596
$ws.Range("K107").Value = 3011
$ws.Range("L107").Value = 3596
$ws.Range("M107").Value = -1091
$ws.Range("N107").Value = -7436

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2824.48
$ws.Range("I31").Value = 1346.3529
$ws.Range("K31").Value = 1346.3529
$ws.Range("M31").Value = -1051.3529
# Row 34
$ws.Range("H34").Value = 2824.48
$ws.Range("I34").Value = 1346.3529
$ws.Range("K34").Value = 1346.3529
$ws.Range("M34").Value = -1144.3529

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 795.3131
$ws.Range("J131").Value = 815.1158
$ws.Range("L131").Value = 2445.3474
$ws.Range("N131").Value = -12525.3474

$ws = $wb.Worksheets.Item("GSM")
# Row 29
$ws.Range("H29").Value = 5096
$ws.Range("J29").Value = 1890
$ws.Range("L29").Value = 1890
$ws.Range("N29").Value = -2470
# Row 70
$ws.Range("H70").Value = 5499.7544
$ws.Range("I70").Value = 5119.425
$ws.Range("J70").Value = 6394.647
$ws.Range("K70").Value = 5119.425
$ws.Range("L70").Value = 6394.647
$ws.Range("M70").Value = -4849.425
$ws.Range("N70").Value = -6934.647
# Row 73
$ws.Range("H73").Value = 5499.7544
$ws.Range("I73").Value = 5119.425
$ws.Range("J73").Value = 6394.647
$ws.Range("K73").Value = 5119.425
$ws.Range("L73").Value = 6394.647
$ws.Range("M73").Value = -4183.425
$ws.Range("N73").Value = -8266.647000000001
# Row 107
$ws.Range("H107").Value = 7408143.5
$ws.Range("I107").Value = 570
$ws.Range("J107").Value = 10101807
$ws.Range("K107").Value = 570
$ws.Range("L107").Value = 10101807
$ws.Range("M107").Value = 1350
$ws.Range("N107").Value = -10105647
# Row 126
$ws.Range("H126").Value = 3264.15
$ws.Range("I126").Value = 2850.6624
$ws.Range("J126").Value = 4648.4346
$ws.Range("K126").Value = 8551.9872
$ws.Range("L126").Value = 13945.3038
$ws.Range("M126").Value = -6081.9872
$ws.Range("N126").Value = -18885.3038
# Row 135
$ws.Range("H135").Value = 21846.154
$ws.Range("J135").Value = 21846.154
$ws.Range("L135").Value = 21846.154
$ws.Range("N135").Value = -31986.154
# Row 138
$ws.Range("H138").Value = 42313.332
$ws.Range("J138").Value = 42313.332
$ws.Range("L138").Value = 42313.332
$ws.Range("N138").Value = -52593.332

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
# Row 100
$ws.Range("H100").Value = 2650.5715
$ws.Range("I100").Value = 1240
$ws.Range("J100").Value = 3434.2222
$ws.Range("K100").Value = 1240
$ws.Range("L100").Value = 3434.2222
$ws.Range("M100").Value = -699
$ws.Range("N100").Value = -4516.2222
# Row 136
$ws.Range("H136").Value = 5804.2
$ws.Range("I136").Value = 1127.1666
$ws.Range("K136").Value = 3381.4998
$ws.Range("M136").Value = -831.4998000000001
# Row 138
$ws.Range("H138").Value = 57990
$ws.Range("J138").Value = 57990
$ws.Range("L138").Value = 57990
$ws.Range("N138").Value = -68270

$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 28266.6
$ws.Range("J15").Value = 28266.6
$ws.Range("L15").Value = 28266.6
$ws.Range("N15").Value = -28842.6
# Row 132
$ws.Range("H132").Value = 7940706
$ws.Range("I132").Value = 4180.8066
$ws.Range("J132").Value = 30307278
$ws.Range("K132").Value = 12542.4198
$ws.Range("L132").Value = 90921834
$ws.Range("M132").Value = -10012.4198
$ws.Range("N132").Value = -90926894
# Row 138
$ws.Range("H138").Value = 41149.668
$ws.Range("J138").Value = 41149.668
$ws.Range("L138").Value = 41149.668
$ws.Range("N138").Value = -51429.668
